$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the (only) worksheet from "Sheet2" to "Sheet5"
$ws.Name = "Sheet5"

# The induction-hardening process now also covers bearing surface 2,
# and its impact (and therefore the overall total) is reduced.
$ws.Range("V16").Value = "Induction Hardening Bearing Surfaces 1, 2"
$ws.Range("V17").Value = 34.074702596165
$ws.Range("C17").Value = 173.688803990365

# Timestamp of the (re-)analysis moved forward by one day / later time of day.
$ws.Range("D1").Value = 45572
$ws.Range("F1").Value = 0.809568472222222
